# Sheet 5 "其他有價證券" (other securities / otherbonds) gets restructured:
#  - header row (row 1) becomes a proper column-header row (name/owner/
#    quantity/face_value/currency/total/property_category/category/date/
#    legislator_name/legislator_id/source_file/index), matching the other
#    sheets' schema (14 columns, same layout as the "股票" sheet).
#  - the data row (row 2) keeps its original values but the quantity (D2)
#    is corrected from 1192.81 to 1192, and gains the same trailing
#    metadata columns (H:N) the other sheets already carry.
#  - the old leftover placeholder row (row 3, which only ever held the
#    stray 財產/種/類/項/件 labels) is removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# --- new header row ------------------------------------------------------
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "owner"
$ws.Range("D1").Value = "quantity"
$ws.Range("E1").Value = "face_value"
$ws.Range("F1").Value = "currency"
$ws.Range("G1").Value = "total"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- fix + extend the data row (row 2) -----------------------------------
$ws.Range("D2").Value = 1192
$ws.Range("H2").Value = "otherbonds"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-11-08"
$ws.Range("K2").Value = "葉宜津"
$ws.Range("L2").Value = 855
$ws.Range("M2").Value = "tmpabd41"
$ws.Range("N2").Value = 109

# --- drop the old leftover placeholder row (old row 3) --------------------
$ws.Rows.Item(3).Delete() | Out-Null
